$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Assets sheet: add Acme/Sha1 asset rows, then trim trailing blank rows ---
$wsAssets.Range("B2").Value = "Acme_UipathURL"
$wsAssets.Range("B3").Value = "Sha1online_Site"
$wsAssets.Range("A2").Value = "Acme_URL"
$wsAssets.Range("A3").Value = "Sha1_URL"

$wsAssets.Rows("4:1000").Delete()

$wsAssets.Range("A4:XFD1048576").Select()

# --- Constants sheet: minor view tweaks ---
$wsConstants.Columns.Item(1).ColumnWidth = 35
$wsConstants.Range("A19:B28").Select()

# --- Settings sheet: add App Credential / App To Kill config rows ---
$wsSettings.Range("B6").Value = "Acme_Cred"
$wsSettings.Range("B6").WrapText = $true
$wsSettings.Range("A6").Value = "Acme_Credential"

$wsSettings.Range("B7").Value = "excel,chrome"
$wsSettings.Range("A7").Value = "App_To_Kill"

$wsSettings.Range("B7").Select()

$wb.Save()
